$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "Rep" column (column D). Remaining columns D (old Date) and
# E (old SO4 amount) shift left automatically, and the now-orphaned
# "Rep"/"stock" shared strings are pruned.
$ws.Columns("D").Delete()

# Fill in the Flask (B) and Dose (C) values that were previously blank
# for most samples.
$flaskDose = @(
    @(5, 0),
    @(26, 5),
    @(35, 7),
    @(24, 3),
    @(18, 1),
    @(1, 0),
    @(17, 1),
    @(24, 5),
    @(28, 5),
    @(21, 3),
    @(27, 5),
    @(19, 5),
    @(8, 0.5),
    @(36, 7),
    @(23, 3),
    @(9, 0.5),
    @(25, 5),
    @(14, 1),
    @(10, 0.5),
    @(3, 0),
    @(31, 7),
    @(33, 7),
    @(34, 7),
    @(16, 1),
    @(30, 5),
    @(22, 3),
    @(7, 0.5),
    @(4, 0),
    @(11, 0.5),
    @(15, 1),
    @(6, 0),
    @(13, 1),
    @(12, 0.5),
    @(20, 3),
    @(32, 7)
)

for ($i = 0; $i -lt $flaskDose.Length; $i++) {
    $row = 2 + $i
    $ws.Range("B$row").Value = $flaskDose[$i][0]
    $ws.Range("C$row").Value = $flaskDose[$i][1]
}

# Row 37 ("Sulfate Stock") gets a label in column B instead of the old
# "stock" note, replacing the text that used to live in column D.
$ws.Range("B37").Value = "sulfate_stock"

# Update the view state (scroll position, zoom, active selection).
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.Zoom = 186
$ws.Range("D41").Select()
